## Generate Report for Handback
## Appends three newly handed-back files to the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"
$include = "Include"

# ---- New files being reported on this run ----
# 1) ffffffc6386561-dede-4fab-89de-d2db3489c788  (shares the ffb66c72 handoff package)
# 2) 6685b6bf-1f52-4832-87df-291ee63b83d0
# 3) bc478b3f-5523-40df-8c0e-efab388c8c5f

$files = @(
    @{
        Name        = "ffffffc6386561-dede-4fab-89de-d2db3489c788"
        ZhXlf       = "ffb66c72-7a23-47c6-82db-d0759d329a7f.41d9e9e7e2c2c44afa6210287a70fd273598c3bf.zh-cn.xlf"
        DeXlf       = "ffb66c72-7a23-47c6-82db-d0759d329a7f.41d9e9e7e2c2c44afa6210287a70fd273598c3bf.de-de.xlf"
        ZhHandoffDt = "2016-01-25 08:35:43"
        ZhHandbackDt= "2016-01-25 08:36:31"
        DeHandoffDt = "2016-01-25 08:35:56"
        DeHandbackDt= "2016-01-25 08:36:52"
        ZhTargetDisplay = "ffb66c72-7a23-47c6-82db-d0759d329a7f.md"
        DeTargetDisplay = "ffb66c72-7a23-47c6-82db-d0759d329a7f.md"
    },
    @{
        Name        = "6685b6bf-1f52-4832-87df-291ee63b83d0"
        ZhXlf       = "6685b6bf-1f52-4832-87df-291ee63b83d0.86fa7517248cbe8736fda64f533993182afad7b5.zh-cn.xlf"
        DeXlf       = "6685b6bf-1f52-4832-87df-291ee63b83d0.86fa7517248cbe8736fda64f533993182afad7b5.de-de.xlf"
        ZhHandoffDt = "2016-01-25 08:39:51"
        ZhHandbackDt= "2016-01-25 08:40:51"
        DeHandoffDt = "2016-01-25 08:40:09"
        DeHandbackDt= "2016-01-25 08:41:14"
        ZhTargetDisplay = "6685b6bf-1f52-4832-87df-291ee63b83d0.md"
        DeTargetDisplay = "6685b6bf-1f52-4832-87df-291ee63b83d0.md"
    },
    @{
        Name        = "bc478b3f-5523-40df-8c0e-efab388c8c5f"
        ZhXlf       = "bc478b3f-5523-40df-8c0e-efab388c8c5f.64ea431d83372bb592cc1fde2022869eebde10aa.zh-cn.xlf"
        DeXlf       = "bc478b3f-5523-40df-8c0e-efab388c8c5f.64ea431d83372bb592cc1fde2022869eebde10aa.de-de.xlf"
        ZhHandoffDt = "2016-01-25 08:40:51"
        ZhHandbackDt= "2016-01-25 08:40:51"
        DeHandoffDt = "2016-01-25 08:40:09"
        DeHandbackDt= "2016-01-25 08:41:14"
        ZhTargetDisplay = "bc478b3f-5523-40df-8c0e-efab388c8c5f.md"
        DeTargetDisplay = "bc478b3f-5523-40df-8c0e-efab388c8c5f.md"
    }
)

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$overviewRepoBase = "https://github.com/OpenLocalizationTest/oltest/blob/451ce9e589f699dbd6bb724e10671ebd88ff894c/e2e"
$zhHandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8cc5869b13bb3c1cf0daa94d4c7cb94e2f686339/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho"
$zhTargetBase = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/16762b7d65c2ced438c1e40ef4ad30791047ba06/e2e"
$zhHandbackBase = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/35a285fbebcd4d6915fadab7b99334576000a875/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho"
$deHandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d988695e908499e2428d8a3b5dd47463b9eba9d6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho"
$deTargetBase = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/bbd46f996869e48f7ef37abe54af1e82be90f090/e2e"
$deHandbackBase = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3488cff47e0ca8c2bf0f68a26387a0ecf76ab925/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho"

$overviewRow = 4
$dataRow = 4

foreach ($f in $files) {

    $mdName = $f.Name + ".md"

    # ---------------- Overview sheet : File Name | zh-cn | de-de ----------------
    $wsOverview.Hyperlinks.Add($wsOverview.Cells.Item($overviewRow, 1), "$overviewRepoBase/$mdName", "", "", $mdName) | Out-Null
    $wsOverview.Cells.Item($overviewRow, 2).Value = $status
    $wsOverview.Cells.Item($overviewRow, 3).Value = $status

    # ---------------- zh-cn sheet ----------------
    # A: Source File Name (hyperlink to source .md)
    $wsZh.Hyperlinks.Add($wsZh.Cells.Item($dataRow, 1), "$overviewRepoBase/$mdName", "", "", $mdName) | Out-Null
    # B: Status
    $wsZh.Cells.Item($dataRow, 2).Value = $status
    # C: Correspond Handoff File (hyperlink to .xlf)
    $wsZh.Hyperlinks.Add($wsZh.Cells.Item($dataRow, 3), "$zhHandoffBase/$($f.ZhXlf)", "", "", $f.ZhXlf) | Out-Null
    # D: Correspond Handoff Datetime
    $wsZh.Cells.Item($dataRow, 4).Value = $f.ZhHandoffDt
    $wsZh.Cells.Item($dataRow, 4).NumberFormat = "yyyy-mm-dd HH:mm:ss"
    # E: Target File (hyperlink to target .md)
    $wsZh.Hyperlinks.Add($wsZh.Cells.Item($dataRow, 5), "$zhTargetBase/$($f.ZhTargetDisplay)", "", "", $f.ZhTargetDisplay) | Out-Null
    # F: Correspond Handback File (hyperlink to .xlf)
    $wsZh.Hyperlinks.Add($wsZh.Cells.Item($dataRow, 6), "$zhHandbackBase/$($f.ZhXlf)", "", "", $f.ZhXlf) | Out-Null
    # G: Correspond Handback DateTime
    $wsZh.Cells.Item($dataRow, 7).Value = $f.ZhHandbackDt
    # H: Handoff Reason
    $wsZh.Cells.Item($dataRow, 8).Value = $include

    # ---------------- de-de sheet ----------------
    # A: Source File Name (hyperlink to source .md)
    $wsDe.Hyperlinks.Add($wsDe.Cells.Item($dataRow, 1), "$overviewRepoBase/$mdName", "", "", $mdName) | Out-Null
    # B: Status
    $wsDe.Cells.Item($dataRow, 2).Value = $status
    # C: Correspond Handoff File (hyperlink to .xlf)
    $wsDe.Hyperlinks.Add($wsDe.Cells.Item($dataRow, 3), "$deHandoffBase/$($f.DeXlf)", "", "", $f.DeXlf) | Out-Null
    # D: Correspond Handoff Datetime
    $wsDe.Cells.Item($dataRow, 4).Value = $f.DeHandoffDt
    $wsDe.Cells.Item($dataRow, 4).NumberFormat = "yyyy-mm-dd HH:mm:ss"
    # E: Target File (hyperlink to target .md)
    $wsDe.Hyperlinks.Add($wsDe.Cells.Item($dataRow, 5), "$deTargetBase/$($f.DeTargetDisplay)", "", "", $f.DeTargetDisplay) | Out-Null
    # F: Correspond Handback File (hyperlink to .xlf)
    $wsDe.Hyperlinks.Add($wsDe.Cells.Item($dataRow, 6), "$deHandbackBase/$($f.DeXlf)", "", "", $f.DeXlf) | Out-Null
    # G: Correspond Handback DateTime
    $wsDe.Cells.Item($dataRow, 7).Value = $f.DeHandbackDt
    # H: Handoff Reason
    $wsDe.Cells.Item($dataRow, 8).Value = $include

    $overviewRow++
    $dataRow++
}

Write-Output "Handback report updated: added $($files.Count) rows to Overview, zh-cn, de-de."
